$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (stored as an Excel date serial number:
# 45415 = 2024-05-03  ->  45436 = 2024-05-24)
$ws.Range("A1").Value2 = 45436

# Update the unit prices in column D (rows 30-33)
$ws.Range("D30").Value2 = 1956.522
$ws.Range("D31").Value2 = 2316.776
$ws.Range("D32").Value2 = 2616.068
$ws.Range("D33").Value2 = 3297.817
